$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Final Test" row (row 25) ready for the user test pass.
$ws.Range("A25").Value = "Final Test"
$ws.Range("B25").Value = 1

# Fix the typo in the if-statement code sample on the Interaction/mouseX task
# (row 17, column E): "false,/b>" -> "false</b>"
$ws.Range("E17").Value = "if(<b style=`"color:#36047c`">false</b>){ <br>
&nbsp;&nbsp;x += random(-1, 1);<br>
&nbsp;&nbsp;y += random(-1, 1);<br>
}"

$ws.Range("C25").Value = "Congratulations on your completion on all the tutorials!<br>
Now Let's use all we've learned to create a interesting interaction effect!
<a href=`"test.html`" target=`"_blank`">Click herel</a> to preview the final effects!
"

$ws.Rows.Item(25).RowHeight = 180

# Move the view/selection down to the newly added row, matching the
# author's final cursor position when the sheet was saved.
$ws.Range("D25").Select()
